$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# BLEU score (row 11)
$ws.Range("B11").Value = 0.05135085245750105

# Code BLEU (row 12)
$ws.Range("B12").Value = 0.1979049916070554
$ws.Range("C12").Value = "{'codebleu': 0.19790499160705538, 'ngram_match_score': 0.05125257286362113, 'weighted_ngram_match_score': 0.09615407512260674, 'syntax_match_score': 0.4519056261343013, 'dataflow_match_score': 0.19230769230769232}"

# Embeddings and Cosine similarity (row 13)
$ws.Range("B13").Value = 0.7952258130289469
